$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update "Semestre ideal" value from EA-8 to EA-7
$ws.Range("B9").Value = "EA-7"
$ws.Range("C9").Value = "EA-7"

# Remove the row containing the "LOB1240 - Condicionantes Geológico" requisite (row 25)
$ws.Rows.Item(25).Delete()
